$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reclassify AFP (row 34) and LDH (row 35) from "hormones" to "pathology"
$ws.Range("G34").Value = "pathology"
$ws.Range("G35").Value = "pathology"

# Add new derived "strata" variables for LDH, AFP and HCG
$ws.Range("A46").Value = "LDH_class"
$ws.Range("B46").Value = "Lactate dehydrogenase, blood concentration strata"
$ws.Range("C46").Value = "LDH strata"
$ws.Range("D46").Value = "LDH strata"
$ws.Range("F46").Value = "factor"
$ws.Range("G46").Value = "pathology"

$ws.Range("A47").Value = "AFP_class"
$ws.Range("B47").Value = "Alpha fetoprotein, blood concentration strata"
$ws.Range("C47").Value = "AFP strata"
$ws.Range("D47").Value = "AFP strata"
$ws.Range("F47").Value = "factor"
$ws.Range("G47").Value = "pathology"

$ws.Range("A48").Value = "HCG_class"
$ws.Range("B48").Value = "Human chorionic gonadotropin, blood concentration strata"
$ws.Range("C48").Value = "HCG strata"
$ws.Range("D48").Value = "HCG strata"
$ws.Range("F48").Value = "factor"
$ws.Range("G48").Value = "hormones"

$ws.Range("G35").Select()
